$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Parts and Trailer Catalog scenarios": reword the part-name text in A4
# and give that cell a wrapped-text style so the longer description displays
# nicely.
$ws.Range("A4").Value = "Replacement Tarp - 18 oz. Super Duty"
$ws.Range("A4").WrapText = $true

# Move the active selection to B4 (as recorded in the saved view state).
$ws.Range("B4").Select() | Out-Null
